$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("D2").Value = "35.409.03"
$ws.Range("E2").Value = "  +0.34%  "
$ws.Range("D3").Value = "1.923.33"
$ws.Range("E3").Value = "  +1.43%  "
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("D5").Value = "'0.737"
$ws.Range("E5").Value = "  +12.97%  "
$ws.Range("D6").Value = "'254.02"
$ws.Range("E6").Value = "  +4.40%  "
$ws.Range("E7").Value = "  -0.15%  "
$ws.Range("D8").Value = "'40.62"
$ws.Range("E8").Value = "  -2.08%  "
$ws.Range("E9").Value = "  +4.28%  "
$ws.Range("E10").Value = "  +5.14%  "
$ws.Range("E11").Value = "  +4.35%  "
$ws.Range("D12").Value = "'0.0998"
$ws.Range("E12").Value = "  +0.05%  "
$ws.Range("E13").Value = "  +1.27%  "
$ws.Range("D14").Value = "'12.70"
$ws.Range("E14").Value = "  +5.44%  "
$ws.Range("E15").Value = "  +3.54%  "
$ws.Range("D16").Value = "1.927.69"
$ws.Range("E16").Value = "  +1.73%  "
$ws.Range("E17").Value = "  +1.75%  "
$ws.Range("D18").Value = "35.439.12"
$ws.Range("E18").Value = "  +0.38%  "
$ws.Range("D19").Value = "'73.58"
$ws.Range("E19").Value = "  +3.21%  "
$ws.Range("D20").Value = "0.0₃0836"
$ws.Range("E20").Value = "  +2.71%  "
$ws.Range("D21").Value = "'13.03"
$ws.Range("E21").Value = "  +4.34%  "
$ws.Range("D22").Value = "'241.80"
$ws.Range("E22").Value = "  +0.21%  "
$ws.Range("D23").Value = "'5.09"
$ws.Range("E23").Value = "  +8.04%  "
$ws.Range("E24").Value = "  -0.17%  "
$ws.Range("E25").Value = "  +2.36%  "
$ws.Range("D26").Value = "'2.36"
$ws.Range("E26").Value = "  -0.74%  "
$ws.Range("D27").Value = "'167.88"
$ws.Range("E27").Value = "  -1.28%  "
$ws.Range("D28").Value = "'8.67"
$ws.Range("E28").Value = "  +4.10%  "
$ws.Range("E29").Value = "  +7.78%  "
$ws.Range("D30").Value = "'18.86"
$ws.Range("E30").Value = "  +3.67%  "
$ws.Range("D31").Value = "4.125.77"
$ws.Range("E31").Value = "  +19.37%  "
$ws.Range("D32").Value = "'4.34"
$ws.Range("E32").Value = "  +5.93%  "
$ws.Range("B33").Value = "WEMIXToken"
$ws.Range("C33").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D33").Value = "'1.99"
$ws.Range("E33").Value = "  +15.00%  "
$ws.Range("B34").Value = "TrustWalletToken"
$ws.Range("C34").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D34").Value = "'1.64"
$ws.Range("E34").Value = "  +23.77%  "
$ws.Range("E35").Value = "  +3.70%  "
$ws.Range("D36").Value = "'4.28"
$ws.Range("E36").Value = "  +4.64%  "
$ws.Range("E37").Value = "  -0.12%  "
$ws.Range("D38").Value = "'0.911"
$ws.Range("E38").Value = "  -1.35%  "
$ws.Range("E39").Value = "  +0.73%  "
$ws.Range("D40").Value = "'17.38"
$ws.Range("E40").Value = "  +10.35%  "
$ws.Range("D41").Value = "'98.94"
$ws.Range("E41").Value = "  +11.06%  "
$ws.Range("E42").Value = "  +4.68%  "
$ws.Range("E43").Value = "  +0.73%  "
$ws.Range("D44").Value = "'0.0651"
$ws.Range("E44").Value = "  +2.52%  "
$ws.Range("E45").Value = "  +5.74%  "
$ws.Range("D46").Value = "1.348.52"
$ws.Range("E46").Value = "  +0.83%  "
$ws.Range("E47").Value = "  +0.81%  "
$ws.Range("B48").Value = "MXToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D48").Value = "'2.77"
$ws.Range("E48").Value = "  +0.01%  "
$ws.Range("B49").Value = "FraxShare"
$ws.Range("C49").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D49").Value = "'6.70"
$ws.Range("E49").Value = "  +3.26%  "
$ws.Range("D50").Value = "'45.39"
$ws.Range("E50").Value = "  -3.15%  "
$ws.Range("D51").Value = "2.109.83"
$ws.Range("E51").Value = "  +1.29%  "
